# Auto-generated edit script applying the cryptos.xlsx update
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "47.319.94"
$ws.Range("E2").Value = "  +2.52%  "

$ws.Range("D3").Value = "2.506.04"
$ws.Range("E3").Value = "  +2.34%  "

$ws.Range("E4").Value = "  +0.06%  "

$ws.Range("D5").Value = "'324.02"
$ws.Range("E5").Value = "  +0.76%  "

$ws.Range("D6").Value = "'109.40"
$ws.Range("E6").Value = "  +4.60%  "

$ws.Range("E7").Value = "  +1.54%  "

$ws.Range("E8").Value = "  +0.02%  "

$ws.Range("E9").Value = "  +0.29%  "

$ws.Range("D10").Value = "'39.24"
$ws.Range("E10").Value = "  +9.21%  "

$ws.Range("E11").Value = "  +1.07%  "

$ws.Range("E12").Value = "  +0.87%  "

$ws.Range("D13").Value = "'18.43"
$ws.Range("E13").Value = "  +0.93%  "

$ws.Range("D14").Value = "'7.21"
$ws.Range("E14").Value = "  +2.18%  "

$ws.Range("D15").Value = "2.897.65"
$ws.Range("E15").Value = "  +2.31%  "

$ws.Range("D16").Value = "2.502.46"
$ws.Range("E16").Value = "  +1.98%  "

$ws.Range("E17").Value = "  +1.90%  "

$ws.Range("D18").Value = "47.252.01"
$ws.Range("E18").Value = "  +2.72%  "

$ws.Range("D19").Value = "'12.87"
$ws.Range("E19").Value = "  +2.05%  "

$ws.Range("E20").Value = "  +4.08%  "

$ws.Range("D21").Value = "0.0₃0943"
$ws.Range("E21").Value = "  +1.25%  "

$ws.Range("E22").Value = "  +14.14%  "

$ws.Range("D23").Value = "'70.52"
$ws.Range("E23").Value = "  -0.50%  "

$ws.Range("D24").Value = "'248.01"
$ws.Range("E24").Value = "  +0.52%  "

$ws.Range("D25").Value = "'2.60"
$ws.Range("E25").Value = "  +3.71%  "

$ws.Range("D26").Value = "'26.04"
$ws.Range("E26").Value = "  +0.66%  "

$ws.Range("D27").Value = "'0.999"
$ws.Range("E27").Value = "  -0.02%  "

$ws.Range("D28").Value = "'2.30"
$ws.Range("E28").Value = "  +0.36%  "

$ws.Range("D29").Value = "'10.05"
$ws.Range("E29").Value = "  +4.04%  "

$ws.Range("D30").Value = "'35.82"
$ws.Range("E30").Value = "  +5.65%  "

$ws.Range("E31").Value = "  +8.07%  "

$ws.Range("D32").Value = "'49.88"
$ws.Range("E32").Value = "  +1.02%  "

$ws.Range("D33").Value = "'20.03"
$ws.Range("E33").Value = "  +0.97%  "

$ws.Range("D34").Value = "'5.46"
$ws.Range("E34").Value = "  +2.28%  "

$ws.Range("E35").Value = "  +4.07%  "

$ws.Range("E36").Value = "  +0.25%  "

$ws.Range("E37").Value = "  +5.28%  "

$ws.Range("E38").Value = "  +4.51%  "

$ws.Range("E39").Value = "  +1.73%  "

$ws.Range("E40").Value = "  +1.26%  "

$ws.Range("E41").Value = "  +0.24%  "

$ws.Range("D42").Value = "'120.80"
$ws.Range("E42").Value = "  -5.25%  "

$ws.Range("D43").Value = "'21.33"
$ws.Range("E43").Value = "  +2.42%  "

$ws.Range("E44").Value = "  +2.49%  "

$ws.Range("D45").Value = "1.998.82"
$ws.Range("E45").Value = "  +1.77%  "

$ws.Range("E46").Value = "  +4.29%  "

$ws.Range("D47").Value = "'2.05"
$ws.Range("E47").Value = "  -1.00%  "

$ws.Range("B48").Value = "Stacks"
$ws.Range("C48").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D48").Value = "'1.78"
$ws.Range("E48").Value = "  -4.26%  "

$ws.Range("B49").Value = "FraxShare"
$ws.Range("C49").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D49").Value = "'9.10"
$ws.Range("E49").Value = "  -0.02%  "

$ws.Range("D50").Value = "'5.22"
$ws.Range("E50").Value = "  +4.18%  "

$ws.Range("D51").Value = "'56.84"
$ws.Range("E51").Value = "  +4.57%  "
